$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# --- Row 11: plain (non-shared) formulas, speed now scaled by 1.5 ---
$ws.Range("E11").Formula = "=1.5*E2/5"
$ws.Range("E11").ClearFormats()
$ws.Range("F11").Formula = "=1.5*F2/5"

# --- Rows 12-15: shared formula block E12:F15 ---
$ws.Range("E12:F15").Formula = "=1.5*E3/5"

# --- Row 16: E16 becomes its own standalone formula, F16 its own formula ---
$ws.Range("E16").Formula = "=1.5*E7/5"
$ws.Range("F16").Formula = "=1.5*F7/5"

# --- Update the active selection shown when the workbook is reopened ---
$ws.Range("E12").Select()
